# Insert a new data row at row 16 (pushing existing rows 16..101 down to 17..102)
# and populate it with the new record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(16).Insert()

$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44687
$ws.Range("D16").NumberFormat = $ws.Range("D17").NumberFormat
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100112021
$ws.Range("G16").Value = "Ají"
$ws.Range("H16").Value = "Inferno"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 25
$ws.Range("K16").Value = 20000
$ws.Range("L16").Value = 21000
$ws.Range("M16").Value = 20600
$ws.Range("N16").Value = "`$/caja 12 kilos"
$ws.Range("O16").Value = "Región de Arica y Parinacota"
$ws.Range("P16").Value = 1717
$ws.Range("Q16").Value = 12
$ws.Range("R16").Value = "Hortaliza"
